# Scheduled runner update: refresh market-derived profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on each leve-profit sheet.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1881.1428
$ws.Cells.Item(19, 9).Value = 1431.6
$ws.Cells.Item(19, 10).Value = 2130.889
$ws.Cells.Item(19, 11).Value = 1431.6
$ws.Cells.Item(19, 12).Value = 2130.889
$ws.Cells.Item(19, 13).Value = -1256.6
$ws.Cells.Item(19, 14).Value = -2480.889
$ws.Cells.Item(33, 8).Value = 200.18182
$ws.Cells.Item(33, 9).Value = 216.64706
$ws.Cells.Item(33, 10).Value = 144.2
$ws.Cells.Item(33, 11).Value = 216.64706
$ws.Cells.Item(33, 12).Value = 144.2
$ws.Cells.Item(33, 13).Value = 12.35293999999999
$ws.Cells.Item(33, 14).Value = -602.2
$ws.Cells.Item(40, 8).Value = 83337140
$ws.Cells.Item(40, 10).Value = 100003570
$ws.Cells.Item(40, 12).Value = 100003570
$ws.Cells.Item(40, 14).Value = -100003920
$ws.Cells.Item(53, 8).Value = 481.24
$ws.Cells.Item(53, 10).Value = 445.35294
$ws.Cells.Item(53, 12).Value = 445.35294
$ws.Cells.Item(53, 14).Value = -1719.35294
$ws.Cells.Item(99, 8).Value = 4119.25
$ws.Cells.Item(99, 9).Value = 154.5
$ws.Cells.Item(99, 10).Value = 5440.8335
$ws.Cells.Item(99, 11).Value = 463.5
$ws.Cells.Item(99, 12).Value = 16322.5005
$ws.Cells.Item(99, 13).Value = 1034.5
$ws.Cells.Item(99, 14).Value = -19318.5005
$ws.Cells.Item(100, 8).Value = 12058.7
$ws.Cells.Item(100, 10).Value = 12398.556
$ws.Cells.Item(100, 12).Value = 12398.556
$ws.Cells.Item(100, 14).Value = -13480.556
$ws.Cells.Item(132, 8).Value = 4439.316
$ws.Cells.Item(132, 9).Value = 2838.16
$ws.Cells.Item(132, 10).Value = 7518.4614
$ws.Cells.Item(132, 11).Value = 8514.48
$ws.Cells.Item(132, 12).Value = 22555.3842
$ws.Cells.Item(132, 13).Value = -5984.48
$ws.Cells.Item(132, 14).Value = -27615.3842

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2587.3157
$ws.Cells.Item(74, 10).Value = 4435.857
$ws.Cells.Item(74, 12).Value = 4435.857
$ws.Cells.Item(74, 14).Value = -6183.857
$ws.Cells.Item(77, 8).Value = 2587.3157
$ws.Cells.Item(77, 10).Value = 4435.857
$ws.Cells.Item(77, 12).Value = 22179.285
$ws.Cells.Item(77, 14).Value = -30915.285
$ws.Cells.Item(97, 8).Value = 1923.4546
$ws.Cells.Item(97, 9).Value = 1416.8235
$ws.Cells.Item(97, 11).Value = 1416.8235
$ws.Cells.Item(97, 13).Value = -920.8235
$ws.Cells.Item(110, 8).Value = 1354.5264
$ws.Cells.Item(110, 9).Value = 449.26666
$ws.Cells.Item(110, 11).Value = 449.26666
$ws.Cells.Item(110, 13).Value = 1595.73334
$ws.Cells.Item(132, 8).Value = 4544.12
$ws.Cells.Item(132, 9).Value = 3784.8948
$ws.Cells.Item(132, 11).Value = 11354.6844
$ws.Cells.Item(132, 13).Value = -8824.6844

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 5346.1904
$ws.Cells.Item(20, 9).Value = 6890.933
$ws.Cells.Item(20, 10).Value = 1484.3334
$ws.Cells.Item(20, 11).Value = 6890.933
$ws.Cells.Item(20, 12).Value = 1484.3334
$ws.Cells.Item(20, 13).Value = -6643.933
$ws.Cells.Item(20, 14).Value = -1978.3334
$ws.Cells.Item(86, 8).Value = 55403.152
$ws.Cells.Item(86, 9).Value = 85339.875
$ws.Cells.Item(86, 10).Value = 7504.4
$ws.Cells.Item(86, 11).Value = 85339.875
$ws.Cells.Item(86, 12).Value = 7504.4
$ws.Cells.Item(86, 13).Value = -84216.875
$ws.Cells.Item(86, 14).Value = -9750.4
$ws.Cells.Item(89, 8).Value = 55403.152
$ws.Cells.Item(89, 9).Value = 85339.875
$ws.Cells.Item(89, 10).Value = 7504.4
$ws.Cells.Item(89, 11).Value = 426699.375
$ws.Cells.Item(89, 12).Value = 37522
$ws.Cells.Item(89, 13).Value = -421083.375
$ws.Cells.Item(89, 14).Value = -48754
$ws.Cells.Item(94, 8).Value = 2325
$ws.Cells.Item(94, 10).Value = 1828.1666
$ws.Cells.Item(94, 12).Value = 1828.1666
$ws.Cells.Item(94, 14).Value = -2730.1666
$ws.Cells.Item(105, 8).Value = 860871.9
$ws.Cells.Item(105, 9).Value = 1608372.5
$ws.Cells.Item(105, 11).Value = 1608372.5
$ws.Cells.Item(105, 13).Value = -1606625.5
$ws.Cells.Item(134, 8).Value = 2761.28
$ws.Cells.Item(134, 9).Value = 2563.4119
$ws.Cells.Item(134, 10).Value = 3181.75
$ws.Cells.Item(134, 11).Value = 7690.2357
$ws.Cells.Item(134, 12).Value = 9545.25
$ws.Cells.Item(134, 13).Value = -5155.2357
$ws.Cells.Item(134, 14).Value = -14615.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 31253218
$ws.Cells.Item(31, 9).Value = 62503244
$ws.Cells.Item(31, 10).Value = 3190.5
$ws.Cells.Item(31, 11).Value = 62503244
$ws.Cells.Item(31, 12).Value = 3190.5
$ws.Cells.Item(31, 13).Value = -62502949
$ws.Cells.Item(31, 14).Value = -3780.5
$ws.Cells.Item(34, 8).Value = 31253218
$ws.Cells.Item(34, 9).Value = 62503244
$ws.Cells.Item(34, 10).Value = 3190.5
$ws.Cells.Item(34, 11).Value = 62503244
$ws.Cells.Item(34, 12).Value = 3190.5
$ws.Cells.Item(34, 13).Value = -62503042
$ws.Cells.Item(34, 14).Value = -3594.5
$ws.Cells.Item(107, 8).Value = 1277.04
$ws.Cells.Item(107, 9).Value = 807.5909
$ws.Cells.Item(107, 10).Value = 4719.6665
$ws.Cells.Item(107, 11).Value = 807.5909
$ws.Cells.Item(107, 12).Value = 4719.6665
$ws.Cells.Item(107, 13).Value = 1112.4091
$ws.Cells.Item(107, 14).Value = -8559.666499999999
$ws.Cells.Item(132, 8).Value = 3066.0417
$ws.Cells.Item(132, 9).Value = 2977.2222
$ws.Cells.Item(132, 11).Value = 8931.6666
$ws.Cells.Item(132, 13).Value = -6401.6666
$ws.Cells.Item(134, 8).Value = 3030.6
$ws.Cells.Item(134, 9).Value = 2852.3635
$ws.Cells.Item(134, 11).Value = 8557.0905
$ws.Cells.Item(134, 13).Value = -6022.0905

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 9662.5
$ws.Cells.Item(69, 9).Value = 2510
$ws.Cells.Item(69, 10).Value = 16815
$ws.Cells.Item(69, 11).Value = 7530
$ws.Cells.Item(69, 12).Value = 50445
$ws.Cells.Item(69, 13).Value = -6719
$ws.Cells.Item(69, 14).Value = -52067
$ws.Cells.Item(72, 8).Value = 9662.5
$ws.Cells.Item(72, 9).Value = 2510
$ws.Cells.Item(72, 10).Value = 16815
$ws.Cells.Item(72, 11).Value = 22590
$ws.Cells.Item(72, 12).Value = 151335
$ws.Cells.Item(72, 13).Value = -18534
$ws.Cells.Item(72, 14).Value = -159447

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2058.4546
$ws.Cells.Item(80, 9).Value = 1605
$ws.Cells.Item(80, 10).Value = 2852
$ws.Cells.Item(80, 11).Value = 1605
$ws.Cells.Item(80, 12).Value = 2852
$ws.Cells.Item(80, 13).Value = -607
$ws.Cells.Item(80, 14).Value = -4848
$ws.Cells.Item(83, 8).Value = 2058.4546
$ws.Cells.Item(83, 9).Value = 1605
$ws.Cells.Item(83, 10).Value = 2852
$ws.Cells.Item(83, 11).Value = 8025
$ws.Cells.Item(83, 12).Value = 14260
$ws.Cells.Item(83, 13).Value = -3033
$ws.Cells.Item(83, 14).Value = -24244
$ws.Cells.Item(97, 8).Value = 636
$ws.Cells.Item(97, 9).Value = 723.75
$ws.Cells.Item(97, 11).Value = 723.75
$ws.Cells.Item(97, 13).Value = -227.75
$ws.Cells.Item(102, 8).Value = 4175.4165
$ws.Cells.Item(102, 9).Value = 4100.5454
$ws.Cells.Item(102, 11).Value = 4100.5454
$ws.Cells.Item(102, 13).Value = -2478.5454
$ws.Cells.Item(107, 8).Value = 1331.45
$ws.Cells.Item(107, 9).Value = 1319.5294
$ws.Cells.Item(107, 11).Value = 1319.5294
$ws.Cells.Item(107, 13).Value = 600.4706000000001
$ws.Cells.Item(126, 8).Value = 3201.087
$ws.Cells.Item(126, 9).Value = 2838.25
$ws.Cells.Item(126, 10).Value = 4030.4285
$ws.Cells.Item(126, 11).Value = 8514.75
$ws.Cells.Item(126, 12).Value = 12091.2855
$ws.Cells.Item(126, 13).Value = -6044.75
$ws.Cells.Item(126, 14).Value = -17031.2855
$ws.Cells.Item(132, 8).Value = 3267.5
$ws.Cells.Item(132, 9).Value = 3267.5
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 9802.5
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -7272.5
$ws.Cells.Item(132, 14).ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 13891222
$ws.Cells.Item(68, 9).Value = 41666664
$ws.Cells.Item(68, 10).Value = 3501.5
$ws.Cells.Item(68, 11).Value = 41666664
$ws.Cells.Item(68, 12).Value = 3501.5
$ws.Cells.Item(68, 13).Value = -41665915
$ws.Cells.Item(68, 14).Value = -4999.5
$ws.Cells.Item(71, 8).Value = 13891222
$ws.Cells.Item(71, 9).Value = 41666664
$ws.Cells.Item(71, 10).Value = 3501.5
$ws.Cells.Item(71, 11).Value = 208333320
$ws.Cells.Item(71, 12).Value = 17507.5
$ws.Cells.Item(71, 13).Value = -208329576
$ws.Cells.Item(71, 14).Value = -24995.5
$ws.Cells.Item(93, 8).Value = 3708600.8
$ws.Cells.Item(93, 9).Value = 3121.3333
$ws.Cells.Item(93, 10).Value = 18530518
$ws.Cells.Item(93, 11).Value = 3121.3333
$ws.Cells.Item(93, 12).Value = 18530518
$ws.Cells.Item(93, 13).Value = -1873.3333
$ws.Cells.Item(93, 14).Value = -18533014
$ws.Cells.Item(132, 8).Value = 4251.467
$ws.Cells.Item(132, 9).Value = 2675.2222
$ws.Cells.Item(132, 10).Value = 6615.8335
$ws.Cells.Item(132, 11).Value = 8025.6666
$ws.Cells.Item(132, 12).Value = 19847.5005
$ws.Cells.Item(132, 13).Value = -5495.6666
$ws.Cells.Item(132, 14).Value = -24907.5005
$ws.Cells.Item(134, 8).Value = 130899
$ws.Cells.Item(134, 10).Value = 130899
$ws.Cells.Item(134, 12).Value = 130899
$ws.Cells.Item(134, 14).Value = -141039
$ws.Cells.Item(136, 8).Value = 3162.682
$ws.Cells.Item(136, 9).Value = 2826.4736
$ws.Cells.Item(136, 10).Value = 5292
$ws.Cells.Item(136, 11).Value = 8479.4208
$ws.Cells.Item(136, 12).Value = 15876
$ws.Cells.Item(136, 13).Value = -5929.4208
$ws.Cells.Item(136, 14).Value = -20976

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 12565
$ws.Cells.Item(96, 10).Value = 15500
$ws.Cells.Item(96, 12).Value = 15500
$ws.Cells.Item(96, 14).Value = -18246
$ws.Cells.Item(126, 8).Value = 5801.3335
$ws.Cells.Item(126, 9).Value = 5801.3335
$ws.Cells.Item(126, 11).Value = 17404.0005
$ws.Cells.Item(126, 13).Value = -14934.0005
$ws.Cells.Item(132, 8).Value = 3676.476
$ws.Cells.Item(132, 9).Value = 3437.796
$ws.Cells.Item(132, 10).Value = 4511.857
$ws.Cells.Item(132, 11).Value = 10313.388
$ws.Cells.Item(132, 12).Value = 13535.571
$ws.Cells.Item(132, 13).Value = -7783.387999999999
$ws.Cells.Item(132, 14).Value = -18595.571
